$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended after the last existing data row (row 78 -> row 79)
# Column A holds the date as text (matching the existing rows, which store
# dates as plain/inline strings rather than Excel date serials), so force
# a text format before assigning the value to prevent Excel from
# auto-converting the "yyyy-mm-dd" looking string into a date serial
# number, then clear the temporary formatting so the cell keeps the
# workbook's default (unstyled) appearance.
$ws.Range("A79").NumberFormat = "@"
$ws.Range("A79").Value = "2025-11-02"
$ws.Range("A79").ClearFormats()

$ws.Range("B79").Value = 59.29999923706055
$ws.Range("C79").Value = 410
$ws.Range("D79").Value = 317.75
